$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Crear Menu de Opciones" task row (row 9): responsible person
# changes from "yoiner" to "breynner" in both the Responsable (D) and last
# (G) columns, reflecting the new "consultar tarea" function assignment.
$ws.Range("D9").Value = "breynner"
$ws.Range("G9").Value = "breynner"

# Update the active cell selection to reflect where the author left off.
$ws.Range("C9").Select()
